$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = $ws.UsedRange.Rows.Count
$cols = $ws.UsedRange.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $val = $cell.Value()
    if ($val -is [string]) {
      $newVal = $val.Replace("D80", "D86").Replace("D51", "D55").Replace("D64", "D69").Replace("S30", "S31")
      if ($newVal -ne $val) {
        $cell.Value = $newVal
      }
    }
  }
}
